$d = $word.ActiveDocument

# Locate the run of text "VLR_BNF1TT_BNFSLD" (the red field name at the end
# of the "Resultado ->" line) and insert a new run "DR" right after it,
# carrying the same red font color, so the two form separate <w:r> runs.
$rng = $d.Content
$rng.Find.Execute("VLR_BNF1TT_BNFSLD", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null

if ($rng.Find.Found) {
    $insertPoint = $rng.Duplicate
    $insertPoint.Collapse(0)  # wdCollapseEnd = 0 -> collapse to end of found range
    $insertPoint.InsertAfter("DR")
    $insertPoint.Font.Color = 255  # wdColorRed
}
